# Insert a new weekly data row for "Choclo" (Dulce o Americano / Primera,
# Región de Arica y Parinacota) at row 585, pushing the existing rows
# 585-630 down to 586-631.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(585).Insert()

$ws.Range("A585").Value = 6
$ws.Range("B585").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C585").Value = "Metropolitana"
$ws.Range("D585").Value = 44461
$ws.Range("E585").Value = 13
$ws.Range("F585").Value = 100112024
$ws.Range("G585").Value = "Choclo"
$ws.Range("H585").Value = "Dulce o Americano"
$ws.Range("I585").Value = "Primera"
$ws.Range("J585").Value = 430
$ws.Range("K585").Value = 32000
$ws.Range("L585").Value = 34000
$ws.Range("M585").Value = 32930
$ws.Range("N585").Value = "`$/malla 70 unidades"
$ws.Range("O585").Value = "Región de Arica y Parinacota"
$ws.Range("P585").Value = 470
$ws.Range("Q585").Value = 70
$ws.Range("R585").Value = "Hortaliza"
